$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-17 (Resolving-Mac as sending cluster) entirely
$ws.Range("A14:T17").EntireRow.Delete() | Out-Null

# Update numeric values for rows 2-13, columns E through T
# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.008000333333333
$ws.Range("H2").Value = 3.024001
$ws.Range("I2").Value = 0.3525296793986107
$ws.Range("J2").Value = 0.3525296793986107
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.844648666666667
$ws.Range("N2").Value = 17.533946
$ws.Range("O2").Value = 0.3204643139023235
$ws.Range("P2").Value = 0.3204643139023235
$ws.Range("Q2").Value = 5.891407804216223
$ws.Range("R2").Value = 53.022670237946
$ws.Range("S2").Value = 0.1129731818386819
$ws.Range("T2").Value = 0.1129731818386819

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.008000333333333
$ws.Range("H3").Value = 3.024001
$ws.Range("I3").Value = 0.3525296793986107
$ws.Range("J3").Value = 0.3525296793986107
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.370261666666667
$ws.Range("N3").Value = 13.110785
$ws.Range("O3").Value = 0.2396231127748354
$ws.Range("P3").Value = 0.2396231127748355
$ws.Range("Q3").Value = 4.405225216753889
$ws.Range("R3").Value = 39.647026950785
$ws.Range("S3").Value = 0.08447425912300989
$ws.Range("T3").Value = 0.08447425912300988

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.008000333333333
$ws.Range("H4").Value = 3.024001
$ws.Range("I4").Value = 0.3525296793986107
$ws.Range("J4").Value = 0.3525296793986107
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.452372666666666
$ws.Range("N4").Value = 19.357118
$ws.Range("O4").Value = 0.3537860524377295
$ws.Range("P4").Value = 0.3537860524377295
$ws.Range("Q4").Value = 6.503993798790889
$ws.Range("R4").Value = 58.53594418911801
$ws.Range("S4").Value = 0.1247200836415729
$ws.Range("T4").Value = 0.1247200836415729

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.008000333333333
$ws.Range("H5").Value = 3.024001
$ws.Range("I5").Value = 0.3525296793986107
$ws.Range("J5").Value = 0.3525296793986107
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.570781
$ws.Range("N5").Value = 4.712343
$ws.Range("O5").Value = 0.08612652088511148
$ws.Range("P5").Value = 0.0861265208851115
$ws.Range("Q5").Value = 1.583347771593667
$ws.Range("R5").Value = 14.250129944343
$ws.Range("S5").Value = 0.0303621547953461
$ws.Range("T5").Value = 0.0303621547953461

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.537538333333333
$ws.Range("H6").Value = 4.612615
$ws.Range("I6").Value = 0.5377259091975243
$ws.Range("J6").Value = 0.5377259091975243
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.844648666666667
$ws.Range("N6").Value = 17.533946
$ws.Range("O6").Value = 0.3204643139023235
$ws.Range("P6").Value = 0.3204643139023235
$ws.Range("Q6").Value = 8.986371369865555
$ws.Range("R6").Value = 80.87734232879
$ws.Range("S6").Value = 0.1723219645584878
$ws.Range("T6").Value = 0.1723219645584878

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.537538333333333
$ws.Range("H7").Value = 4.612615
$ws.Range("I7").Value = 0.5377259091975243
$ws.Range("J7").Value = 0.5377259091975243
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.370261666666667
$ws.Range("N7").Value = 13.110785
$ws.Range("O7").Value = 0.2396231127748354
$ws.Range("P7").Value = 0.2396231127748355
$ws.Range("Q7").Value = 6.719444839197222
$ws.Range("R7").Value = 60.475003552775
$ws.Range("S7").Value = 0.1288515561815893
$ws.Range("T7").Value = 0.1288515561815893

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.537538333333333
$ws.Range("H8").Value = 4.612615
$ws.Range("I8").Value = 0.5377259091975243
$ws.Range("J8").Value = 0.5377259091975243
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.452372666666666
$ws.Range("N8").Value = 19.357118
$ws.Range("O8").Value = 0.3537860524377295
$ws.Range("P8").Value = 0.3537860524377295
$ws.Range("Q8").Value = 9.920770315952222
$ws.Range("R8").Value = 89.28693284357
$ws.Range("S8").Value = 0.1902399267084811
$ws.Range("T8").Value = 0.1902399267084811

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.537538333333333
$ws.Range("H9").Value = 4.612615
$ws.Range("I9").Value = 0.5377259091975243
$ws.Range("J9").Value = 0.5377259091975243
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.570781
$ws.Range("N9").Value = 4.712343
$ws.Range("O9").Value = 0.08612652088511148
$ws.Range("P9").Value = 0.0861265208851115
$ws.Range("Q9").Value = 2.415136000771667
$ws.Range("R9").Value = 21.736224006945
$ws.Range("S9").Value = 0.04631246174896614
$ws.Range("T9").Value = 0.04631246174896615

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.313796
$ws.Range("H10").Value = 0.941388
$ws.Range("I10").Value = 0.1097444114038651
$ws.Range("J10").Value = 0.1097444114038651
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.844648666666667
$ws.Range("N10").Value = 17.533946
$ws.Range("O10").Value = 0.3204643139023235
$ws.Range("P10").Value = 0.3204643139023235
$ws.Range("Q10").Value = 1.834027373005334
$ws.Range("R10").Value = 16.506246357048
$ws.Range("S10").Value = 0.03516916750515395
$ws.Range("T10").Value = 0.03516916750515395

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.313796
$ws.Range("H11").Value = 0.941388
$ws.Range("I11").Value = 0.1097444114038651
$ws.Range("J11").Value = 0.1097444114038651
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 4.370261666666667
$ws.Range("N11").Value = 13.110785
$ws.Range("O11").Value = 0.2396231127748354
$ws.Range("P11").Value = 0.2396231127748355
$ws.Range("Q11").Value = 1.371370629953333
$ws.Range("R11").Value = 12.34233566958
$ws.Range("S11").Value = 0.0262972974702363
$ws.Range("T11").Value = 0.0262972974702363

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.313796
$ws.Range("H12").Value = 0.941388
$ws.Range("I12").Value = 0.1097444114038651
$ws.Range("J12").Value = 0.1097444114038651
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 6.452372666666666
$ws.Range("N12").Value = 19.357118
$ws.Range("O12").Value = 0.3537860524377295
$ws.Range("P12").Value = 0.3537860524377295
$ws.Range("Q12").Value = 2.024728733309333
$ws.Range("R12").Value = 18.222558599784
$ws.Range("S12").Value = 0.03882604208767557
$ws.Range("T12").Value = 0.03882604208767557

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.313796
$ws.Range("H13").Value = 0.941388
$ws.Range("I13").Value = 0.1097444114038651
$ws.Range("J13").Value = 0.1097444114038651
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.570781
$ws.Range("N13").Value = 4.712343
$ws.Range("O13").Value = 0.08612652088511148
$ws.Range("P13").Value = 0.0861265208851115
$ws.Range("Q13").Value = 0.492904794676
$ws.Range("R13").Value = 4.436143152084
$ws.Range("S13").Value = 0.009451904340799252
$ws.Range("T13").Value = 0.009451904340799252
